# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. Update DAMSLTag (col I) and DialogAct (col J)
# values for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 11; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 16; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 22; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 41; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 44; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 46; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 47; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 51; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 70; Tag = "sd"; Act = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
